# Auto-committed on 2022/03/09 週三
#
# This script reproduces the edits made to CdIndustry.xlsx:
#  - DBS sheet (sheet2): row 3, column B text changed from the old
#    "IndustryCode >= ,AND IndustryCode <= " range-filter string to the new
#    "IndustryCode %" wildcard-filter string.
#  - DBS sheet: a new row 4 added containing a new lookup definition
#    ("findIndustryItem" / "IndustryItem %" / "IndustryCode asc").
#  - The DBS sheet becomes the active/selected sheet (tab), replacing DBD.
#  - Selection markers on both sheets move to reflect where the user was
#    last working (DBD -> B10, DBS -> B8).

$wb = $excel.ActiveWorkbook

$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# --- DBS (sheet2): append new row 4, then update existing row 3 -----------

# New row describing a second lookup function for IndustryItem.
$dbs.Range("A4").Value = "findIndustryItem"
$dbs.Range("B4").Value = "IndustryItem %"
$dbs.Range("C4").Value = "IndustryCode asc"

# B3 previously held the now-unused "IndustryCode >= ,AND IndustryCode <= "
# string; it is replaced with the new wildcard-search label.
$dbs.Range("B3").Value = "IndustryCode %"

# --- Selections / active sheet --------------------------------------------

# Move the DBD selection before switching away from it so its saved
# selection reflects where editing left off.
[void]$dbd.Range("B10").Select()

# DBS becomes the active (selected) sheet/tab.
[void]$dbs.Activate()
[void]$dbs.Range("B8").Select()
